# Update the "want to go" counts (column F) on the "展览" and "全部类型" sheets
# Row 4 (F4): 288 -> 290
# Row 5 (F5): 4236 -> 4252

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F4").Value = 290
    $ws.Range("F5").Value = 4252
}
